# "Adding Master Data XLS"
#
# The existing "DNG / Dongle" master-data rows (eng/ara/fra) are replaced
# with "DKS / Dekstop / Desktop Computer" rows in the same three languages.
# code (A) and lang_code (D) stay put; name (B) and descr (C) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - English
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"
$ws.Range("D2").Value = "eng"

# Row 3 - Arabic
$ws.Range("A3").Value = "DKS"
$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"
$ws.Range("D3").Value = "ara"

# Row 4 - French
$ws.Range("A4").Value = "DKS"
$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"
$ws.Range("D4").Value = "fra"

# Print setup (A4 -> Letter/portrait, as saved from the author's machine)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor parked on D10, matching the saved selection
$ws.Range("D10").Select() | Out-Null
